# Auto-generated edit script: updates cryptos price (D) and volume (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): force-write as text (NumberFormat @) so purely-numeric-looking
# strings (e.g. '229.01') are not silently reinterpreted as numbers, then restore
# the default 'Normal' style so no stray formatting is left on the cell.
$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '38.762.65'
$c.Style = "Normal"
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.093.14'
$c.Style = "Normal"
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '229.01'
$c.Style = "Normal"
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.617'
$c.Style = "Normal"
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '61.13'
$c.Style = "Normal"
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.386'
$c.Style = "Normal"
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.0844'
$c.Style = "Normal"
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '15.31'
$c.Style = "Normal"
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '2.404.10'
$c.Style = "Normal"
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '22.01'
$c.Style = "Normal"
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.804'
$c.Style = "Normal"
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '5.49'
$c.Style = "Normal"
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '2.094.06'
$c.Style = "Normal"
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '38.705.28'
$c.Style = "Normal"
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '71.75'
$c.Style = "Normal"
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '6.08'
$c.Style = "Normal"
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '0.0₃0842'
$c.Style = "Normal"
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '227.82'
$c.Style = "Normal"
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.34'
$c.Style = "Normal"
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '171.45'
$c.Style = "Normal"
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '9.53'
$c.Style = "Normal"
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '0.137'
$c.Style = "Normal"
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '1.42'
$c.Style = "Normal"
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '19.32'
$c.Style = "Normal"
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '4.52'
$c.Style = "Normal"
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '4.74'
$c.Style = "Normal"
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.0613'
$c.Style = "Normal"
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.0228'
$c.Style = "Normal"
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '100.94'
$c.Style = "Normal"
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '1.535.46'
$c.Style = "Normal"
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.0910'
$c.Style = "Normal"
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '7.67'
$c.Style = "Normal"
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '2.96'
$c.Style = "Normal"
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '2.291.21'
$c.Style = "Normal"

# Volume(1h) column (E): plain text values (never numeric-looking), set directly.
$ws.Range('E2').Value = '  +1.33%  '
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +1.38%  '
$ws.Range('E10').Value = '  +0.21%  '
$ws.Range('E12').Value = '  +4.08%  '
$ws.Range('E13').Value = '  -0.23%  '
$ws.Range('E14').Value = '  -2.09%  '
$ws.Range('E15').Value = '  +3.52%  '
$ws.Range('E16').Value = '  -0.28%  '
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('E18').Value = '  +1.46%  '
$ws.Range('E19').Value = '  +1.99%  '
$ws.Range('E20').Value = '  +0.97%  '
$ws.Range('E21').Value = '  +0.52%  '
$ws.Range('E22').Value = '  +1.46%  '
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('E27').Value = '  +0.83%  '
$ws.Range('E28').Value = '  +4.33%  '
$ws.Range('E29').Value = '  +5.86%  '
$ws.Range('E30').Value = '  +1.22%  '
$ws.Range('E31').Value = '  +2.73%  '
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('E33').Value = '  +1.75%  '
$ws.Range('E34').Value = '  +0.82%  '
$ws.Range('E35').Value = '  +1.02%  '
$ws.Range('E36').Value = '  -1.02%  '
$ws.Range('E37').Value = '  -0.55%  '
$ws.Range('E38').Value = '  +1.22%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('E41').Value = '  +4.17%  '
$ws.Range('E42').Value = '  +0.82%  '
$ws.Range('E43').Value = '  -0.80%  '
$ws.Range('E44').Value = '  -1.01%  '
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('E46').Value = '  +1.38%  '
$ws.Range('E47').Value = '  +5.49%  '
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('E50').Value = '  -1.48%  '
$ws.Range('E51').Value = '  -0.21%  '
